$wb = $excel.ActiveWorkbook

# Append "Sheet3" after the last existing sheet; A1 is touched but left empty
$lastIndex = $wb.Worksheets.Count
$sheet3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$sheet3.Name = "Sheet3"
$sheet3.Range("A1").Font.Bold = $false

# Append "Sheet4" after Sheet3; A1 holds the existing shared string "No of Job Seekers"
$lastIndex = $wb.Worksheets.Count
$sheet4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$sheet4.Name = "Sheet4"
$sheet4.Range("A1").Value = "No of Job Seekers"
